$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("H17").Value = 33673.965
$ws.Range("J17").Value = 34949.883
$ws.Range("L17").Value = 104849.649
$ws.Range("N17").Value = -105185.649
$ws.Range("H33").Value = 959251.5
$ws.Range("I33").Value = 1231809.1
$ws.Range("K33").Value = 1231809.1
$ws.Range("M33").Value = -1231580.1
$ws.Range("H38").Value = 572.8
$ws.Range("I38").Value = 466
$ws.Range("K38").Value = 1398
$ws.Range("M38").Value = -1026
$ws.Range("H62").Value = 24550.727
$ws.Range("J62").Value = 24009.666
$ws.Range("L62").Value = 24009.666
$ws.Range("N62").Value = -25257.666
$ws.Range("H65").Value = 24550.727
$ws.Range("J65").Value = 24009.666
$ws.Range("L65").Value = 120048.33
$ws.Range("N65").Value = -126288.33
$ws.Range("H76").Value = 3695.8
$ws.Range("I76").Value = 3994.75
$ws.Range("J76").Value = 2500
$ws.Range("K76").Value = 3994.75
$ws.Range("L76").Value = 2500
$ws.Range("M76").Value = -3679.75
$ws.Range("N76").Value = -3130
$ws.Range("H79").Value = 3695.8
$ws.Range("I79").Value = 3994.75
$ws.Range("J79").Value = 2500
$ws.Range("K79").Value = 3994.75
$ws.Range("L79").Value = 2500
$ws.Range("M79").Value = -2902.75
$ws.Range("N79").Value = -4684
$ws.Range("H80").Value = 1144
$ws.Range("I80").Value = 1149.4
$ws.Range("J80").Value = 1137.25
$ws.Range("K80").Value = 3448.2
$ws.Range("L80").Value = 3411.75
$ws.Range("M80").Value = -2450.2
$ws.Range("N80").Value = -5407.75
$ws.Range("H83").Value = 1144
$ws.Range("I83").Value = 1149.4
$ws.Range("J83").Value = 1137.25
$ws.Range("K83").Value = 10344.6
$ws.Range("L83").Value = 10235.25
$ws.Range("M83").Value = -5352.6
$ws.Range("N83").Value = -20219.25
$ws.Range("H113").Value = 5368.8945
$ws.Range("I113").Value = 5252.25
$ws.Range("J113").Value = 5400
$ws.Range("K113").Value = 5252.25
$ws.Range("L113").Value = 5400
$ws.Range("M113").Value = -1998.25
$ws.Range("N113").Value = -11908
$ws.Range("H116").Value = 30594
$ws.Range("I116").Value = 32712.666
$ws.Range("K116").Value = 32712.666
$ws.Range("M116").Value = -29270.666
$ws.Range("H127").Value = 43453.816
$ws.Range("I127").Value = 52825.945
$ws.Range("J127").Value = 1279.25
$ws.Range("K127").Value = 158477.835
$ws.Range("L127").Value = 3837.75
$ws.Range("M127").Value = -153517.835
$ws.Range("N127").Value = -13757.75
$ws.Range("H132").Value = 2637.7856
$ws.Range("I132").Value = 1159.7949
$ws.Range("J132").Value = 21851.666
$ws.Range("K132").Value = 3479.384700000001
$ws.Range("L132").Value = 65554.99800000001
$ws.Range("M132").Value = -949.3847000000005
$ws.Range("N132").Value = -70614.99800000001
$ws.Range("H135").Value = 515.5625
$ws.Range("I135").Value = 426.84616
$ws.Range("K135").Value = 3841.61544
$ws.Range("M135").Value = -1306.61544
$ws.Range("H138").Value = 2438.574
$ws.Range("I138").Value = 1255.8
$ws.Range("J138").Value = 3134.3235
$ws.Range("K138").Value = 3767.4
$ws.Range("L138").Value = 9402.970499999999
$ws.Range("M138").Value = 1372.6
$ws.Range("N138").Value = -19682.9705
$ws.Range("H141").Value = 38933.875
$ws.Range("I141").Value = 38933.875
$ws.Range("K141").Value = 116801.625
$ws.Range("M141").Value = -111621.625
$ws.Range("M11").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2376.2727
$ws.Range("I2").Value = 2071.5557
$ws.Range("K2").Value = 2071.5557
$ws.Range("M2").Value = -1958.5557
$ws.Range("H21").Value = 1839.25
$ws.Range("I21").Value = 679
$ws.Range("K21").Value = 679
$ws.Range("M21").Value = -305
$ws.Range("H32").Value = 8609783
$ws.Range("I32").Value = 1623507.2
$ws.Range("J32").Value = 31315180
$ws.Range("K32").Value = 1623507.2
$ws.Range("L32").Value = 31315180
$ws.Range("M32").Value = -1623220.2
$ws.Range("N32").Value = -31315754
$ws.Range("H45").Value = 10597.315
$ws.Range("I45").Value = 9788.963
$ws.Range("J45").Value = 12581.454
$ws.Range("K45").Value = 9788.963
$ws.Range("L45").Value = 12581.454
$ws.Range("M45").Value = -9411.963
$ws.Range("N45").Value = -13335.454
$ws.Range("H61").Value = 1550.6072
$ws.Range("I61").Value = 1550.6072
$ws.Range("K61").Value = 1550.6072
$ws.Range("M61").Value = -1338.6072
$ws.Range("H74").Value = 1334.0731
$ws.Range("I74").Value = 1331.5294
$ws.Range("K74").Value = 1331.5294
$ws.Range("M74").Value = -457.5293999999999
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("H77").Value = 1334.0731
$ws.Range("I77").Value = 1331.5294
$ws.Range("K77").Value = 6657.646999999999
$ws.Range("M77").Value = -2289.646999999999
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("H116").Value = 2376.2727
$ws.Range("I116").Value = 2071.5557
$ws.Range("K116").Value = 2071.5557
$ws.Range("M116").Value = 222.4443000000001
$ws.Range("H132").Value = 5339.857
$ws.Range("I132").Value = 5500.905
$ws.Range("J132").Value = 4856.7144
$ws.Range("K132").Value = 16502.715
$ws.Range("L132").Value = 14570.1432
$ws.Range("M132").Value = -13972.715
$ws.Range("N132").Value = -19630.1432
$ws.Range("H136").Value = 1550.6072
$ws.Range("I136").Value = 1550.6072
$ws.Range("K136").Value = 4651.821599999999
$ws.Range("M136").Value = -2101.821599999999
$ws.Range("N76").ClearContents()
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2376.2727
$ws.Range("I3").Value = 2071.5557
$ws.Range("K3").Value = 2071.5557
$ws.Range("M3").Value = -1957.5557
$ws.Range("H92").Value = 40267.332
$ws.Range("J92").Value = 40267.332
$ws.Range("L92").Value = 40267.332
$ws.Range("N92").Value = -45259.332
$ws.Range("H117").Value = 100742
$ws.Range("J117").Value = 100742
$ws.Range("L117").Value = 100742
$ws.Range("N117").Value = -109920
$ws.Range("H134").Value = 3322.5
$ws.Range("I134").Value = 2146.9512
$ws.Range("K134").Value = 6440.8536
$ws.Range("M134").Value = -3905.8536

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 579
$ws.Range("I7").Value = 399.5
$ws.Range("K7").Value = 399.5
$ws.Range("M7").Value = -286.5
$ws.Range("H31").Value = 3832.2
$ws.Range("I31").Value = 2053.3157
$ws.Range("K31").Value = 2053.3157
$ws.Range("M31").Value = -1758.3157
$ws.Range("H32").Value = 6704
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 7644.8
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 7644.8
$ws.Range("M32").Value = -1684
$ws.Range("N32").Value = -8276.799999999999
$ws.Range("H34").Value = 3832.2
$ws.Range("I34").Value = 2053.3157
$ws.Range("K34").Value = 2053.3157
$ws.Range("M34").Value = -1851.3157
$ws.Range("H50").Value = 9444
$ws.Range("J50").Value = 9999.5
$ws.Range("L50").Value = 9999.5
$ws.Range("N50").Value = -11249.5
$ws.Range("H51").Value = 9999.333000000001
$ws.Range("J51").Value = 9999.333000000001
$ws.Range("L51").Value = 9999.333000000001
$ws.Range("N51").Value = -11471.333
$ws.Range("H61").Value = 9999.333000000001
$ws.Range("J61").Value = 9999.333000000001
$ws.Range("L61").Value = 9999.333000000001
$ws.Range("N61").Value = -10695.333
$ws.Range("H107").Value = 1046.5294
$ws.Range("I107").Value = 829.9091
$ws.Range("K107").Value = 829.9091
$ws.Range("M107").Value = 1090.0909
$ws.Range("H132").Value = 2708.7144
$ws.Range("I132").Value = 2538.75
$ws.Range("K132").Value = 7616.25
$ws.Range("M132").Value = -5086.25
$ws.Range("H134").Value = 1978
$ws.Range("I134").Value = 1688.1786
$ws.Range("K134").Value = 5064.5358
$ws.Range("M134").Value = -2529.5358
$ws.Range("H138").Value = 64889.5
$ws.Range("J138").Value = 64889.5
$ws.Range("L138").Value = 64889.5
$ws.Range("N138").Value = -75169.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 700
$ws.Range("I7").Value = 1040
$ws.Range("K7").Value = 3120
$ws.Range("M7").Value = -3008
$ws.Range("H68").Value = 1391031.6
$ws.Range("J68").Value = 2176015.5
$ws.Range("L68").Value = 6528046.5
$ws.Range("N68").Value = -6529668.5
$ws.Range("H71").Value = 1391031.6
$ws.Range("J71").Value = 2176015.5
$ws.Range("L71").Value = 19584139.5
$ws.Range("N71").Value = -19592251.5
$ws.Range("H80").Value = 11780570
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 6000
$ws.Range("M80").Value = -5064
$ws.Range("H83").Value = 11780570
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 18000
$ws.Range("M83").Value = -13320
$ws.Range("H107").Value = 983.5625
$ws.Range("I107").Value = 501.33334
$ws.Range("J107").Value = 1094.8462
$ws.Range("K107").Value = 1504.00002
$ws.Range("L107").Value = 3284.5386
$ws.Range("M107").Value = 415.9999800000001
$ws.Range("N107").Value = -7124.5386
$ws.Range("H131").Value = 82243.44
$ws.Range("I131").Value = 53671.42
$ws.Range("J131").Value = 106919.27
$ws.Range("K131").Value = 161014.26
$ws.Range("L131").Value = 320757.81
$ws.Range("M131").Value = -155974.26
$ws.Range("N131").Value = -330837.81
$ws.Range("H132").Value = 1999.8572
$ws.Range("J132").Value = 1999.6666
$ws.Range("L132").Value = 17996.9994
$ws.Range("N132").Value = -23056.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12217
$ws.Range("I80").Value = 18733.834
$ws.Range("J80").Value = 4396.8
$ws.Range("K80").Value = 18733.834
$ws.Range("L80").Value = 4396.8
$ws.Range("M80").Value = -17735.834
$ws.Range("N80").Value = -6392.8
$ws.Range("H83").Value = 12217
$ws.Range("I83").Value = 18733.834
$ws.Range("J83").Value = 4396.8
$ws.Range("K83").Value = 93669.17
$ws.Range("L83").Value = 21984
$ws.Range("M83").Value = -88677.17
$ws.Range("N83").Value = -31968
$ws.Range("H102").Value = 3264.3914
$ws.Range("I102").Value = 2793.0527
$ws.Range("J102").Value = 5503.25
$ws.Range("K102").Value = 2793.0527
$ws.Range("L102").Value = 5503.25
$ws.Range("M102").Value = -1171.0527
$ws.Range("N102").Value = -8747.25
$ws.Range("H107").Value = 16534.064
$ws.Range("I107").Value = 24596.316
$ws.Range("J107").Value = 3768.8333
$ws.Range("K107").Value = 24596.316
$ws.Range("L107").Value = 3768.8333
$ws.Range("M107").Value = -22676.316
$ws.Range("N107").Value = -7608.8333
$ws.Range("H113").Value = 12180.75
$ws.Range("I113").Value = 17027.428
$ws.Range("J113").Value = 5395.4
$ws.Range("K113").Value = 17027.428
$ws.Range("L113").Value = 5395.4
$ws.Range("M113").Value = -14857.428
$ws.Range("N113").Value = -9735.4
$ws.Range("H132").Value = 4143.865
$ws.Range("I132").Value = 4241.1953
$ws.Range("K132").Value = 12723.5859
$ws.Range("M132").Value = -10193.5859

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2070.1
$ws.Range("I16").Value = 2000.1428
$ws.Range("J16").Value = 2233.3333
$ws.Range("K16").Value = 2000.1428
$ws.Range("L16").Value = 2233.3333
$ws.Range("M16").Value = -1830.1428
$ws.Range("N16").Value = -2573.3333
$ws.Range("H22").Value = 988.38464
$ws.Range("J22").Value = 1900
$ws.Range("L22").Value = 1900
$ws.Range("N22").Value = -2490
$ws.Range("H27").Value = 988.38464
$ws.Range("J27").Value = 1900
$ws.Range("L27").Value = 1900
$ws.Range("N27").Value = -2114
$ws.Range("H46").Value = 958.2143
$ws.Range("I46").Value = 971.2857
$ws.Range("K46").Value = 971.2857
$ws.Range("M46").Value = -783.2857
$ws.Range("H53").Value = 35998.332
$ws.Range("I53").Value = 39000
$ws.Range("K53").Value = 39000
$ws.Range("M53").Value = -38482
$ws.Range("H61").Value = 37042436
$ws.Range("I61").Value = 47623136
$ws.Range("J61").Value = 9998.5
$ws.Range("K61").Value = 47623136
$ws.Range("L61").Value = 9998.5
$ws.Range("M61").Value = -47622934
$ws.Range("N61").Value = -10402.5
$ws.Range("H68").Value = 21115.938
$ws.Range("I68").Value = 2305.8
$ws.Range("J68").Value = 52466.168
$ws.Range("K68").Value = 2305.8
$ws.Range("L68").Value = 52466.168
$ws.Range("M68").Value = -1556.8
$ws.Range("N68").Value = -53964.168
$ws.Range("H71").Value = 21115.938
$ws.Range("I71").Value = 2305.8
$ws.Range("J71").Value = 52466.168
$ws.Range("K71").Value = 11529
$ws.Range("L71").Value = 262330.84
$ws.Range("M71").Value = -7785
$ws.Range("N71").Value = -269818.84
$ws.Range("H100").Value = 108870.164
$ws.Range("I100").Value = 63977.6
$ws.Range("K100").Value = 63977.6
$ws.Range("M100").Value = -63436.6
$ws.Range("H113").Value = 37042436
$ws.Range("I113").Value = 47623136
$ws.Range("J113").Value = 9998.5
$ws.Range("K113").Value = 47623136
$ws.Range("L113").Value = 9998.5
$ws.Range("M113").Value = -47620966
$ws.Range("N113").Value = -14338.5
$ws.Range("H122").Value = 5899
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 21930
$ws.Range("J20").Value = 22702.2
$ws.Range("L20").Value = 22702.2
$ws.Range("N20").Value = -23182.2
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("H53").Value = 3100
$ws.Range("I53").Value = 3100
$ws.Range("K53").Value = 3100
$ws.Range("M53").Value = -2493
$ws.Range("H81").Value = 502500000
$ws.Range("I81").Value = 502500000
$ws.Range("K81").Value = 1005000000
$ws.Range("M81").Value = -1004998939
$ws.Range("H84").Value = 502500000
$ws.Range("I84").Value = 502500000
$ws.Range("K84").Value = 5025000000
$ws.Range("M84").Value = -5024994696
$ws.Range("H126").Value = 993.7059
$ws.Range("I126").Value = 993.7059
$ws.Range("K126").Value = 2981.1177
$ws.Range("M126").Value = -511.1177000000002
$ws.Range("H132").Value = 6711.148
$ws.Range("I132").Value = 7309.625
$ws.Range("J132").Value = 1923.3334
$ws.Range("K132").Value = 21928.875
$ws.Range("L132").Value = 5770.0002
$ws.Range("M132").Value = -19398.875
$ws.Range("N132").Value = -10830.0002
$ws.Range("H133").Value = 65867.8
$ws.Range("J133").Value = 65867.8
$ws.Range("L133").Value = 65867.8
$ws.Range("N133").Value = -75987.8
$ws.Range("H136").Value = 2269.1538
$ws.Range("I136").Value = 2321.6738
$ws.Range("J136").Value = 1866.5
$ws.Range("K136").Value = 6965.0214
$ws.Range("L136").Value = 5599.5
$ws.Range("M136").Value = -4415.0214
$ws.Range("N136").Value = -10699.5
$ws.Range("N30").ClearContents()
